# update scripts wuth new tpm
# Recomputed NATMI ligand-receptor edge statistics (Rln1-Rxfp2) for the
# YoungD0 TPM run: the 3-cluster (ECs/FAPs/MuSCs) x 2-target (FAPs/MuSCs)
# cross product is now written out in full (6 data rows instead of 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 6,20

# Row 2: ECs -> FAPs
$data[0,0]  = "ECs"
$data[0,1]  = "Rln1"
$data[0,2]  = "Rxfp2"
$data[0,3]  = "FAPs"
$data[0,4]  = 3
$data[0,5]  = 1
$data[0,6]  = 0.08104266666666667
$data[0,7]  = 0.243128
$data[0,8]  = 0.1109852216299026
$data[0,9]  = 0.1109852216299026
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.5698856666666666
$data[0,13] = 1.709657
$data[0,14] = 0.9886227745742286
$data[0,15] = 0.9886227745742288
$data[0,16] = 0.04618505412177777
$data[0,17] = 0.415665487096
$data[0,18] = 0.10972251774449
$data[0,19] = 0.10972251774449

# Row 3: ECs -> MuSCs
$data[1,0]  = "ECs"
$data[1,1]  = "Rln1"
$data[1,2]  = "Rxfp2"
$data[1,3]  = "MuSCs"
$data[1,4]  = 3
$data[1,5]  = 1
$data[1,6]  = 0.08104266666666667
$data[1,7]  = 0.243128
$data[1,8]  = 0.1109852216299026
$data[1,9]  = 0.1109852216299026
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.006558333333333333
$data[1,13] = 0.019675
$data[1,14] = 0.01137722542577134
$data[1,15] = 0.01137722542577134
$data[1,16] = 0.0005315048222222222
$data[1,17] = 0.0047835434
$data[1,18] = 0.001262703885412594
$data[1,19] = 0.001262703885412595

# Row 4: FAPs -> FAPs
$data[2,0]  = "FAPs"
$data[2,1]  = "Rln1"
$data[2,2]  = "Rxfp2"
$data[2,3]  = "FAPs"
$data[2,4]  = 3
$data[2,5]  = 1
$data[2,6]  = 0.5478883333333334
$data[2,7]  = 1.643665
$data[2,8]  = 0.7503147490635131
$data[2,9]  = 0.7503147490635131
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.5698856666666666
$data[2,13] = 1.709657
$data[2,14] = 0.9886227745742286
$data[2,15] = 0.9886227745742288
$data[2,16] = 0.3122337081005556
$data[2,17] = 2.810103372905
$data[2,18] = 0.7417782490231365
$data[2,19] = 0.7417782490231366

# Row 5: FAPs -> MuSCs
$data[3,0]  = "FAPs"
$data[3,1]  = "Rln1"
$data[3,2]  = "Rxfp2"
$data[3,3]  = "MuSCs"
$data[3,4]  = 3
$data[3,5]  = 1
$data[3,6]  = 0.5478883333333334
$data[3,7]  = 1.643665
$data[3,8]  = 0.7503147490635131
$data[3,9]  = 0.7503147490635131
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.006558333333333333
$data[3,13] = 0.019675
$data[3,14] = 0.01137722542577134
$data[3,15] = 0.01137722542577134
$data[3,16] = 0.003593234319444445
$data[3,17] = 0.032339108875
$data[3,18] = 0.008536500040376642
$data[3,19] = 0.008536500040376642

# Row 6: MuSCs -> FAPs
$data[4,0]  = "MuSCs"
$data[4,1]  = "Rln1"
$data[4,2]  = "Rxfp2"
$data[4,3]  = "FAPs"
$data[4,4]  = 3
$data[4,5]  = 1
$data[4,6]  = 0.1012803333333333
$data[4,7]  = 0.303841
$data[4,8]  = 0.1387000293065843
$data[4,9]  = 0.1387000293065843
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.5698856666666666
$data[4,13] = 1.709657
$data[4,14] = 0.9886227745742286
$data[4,15] = 0.9886227745742288
$data[4,16] = 0.05771821028188889
$data[4,17] = 0.5194638925370001
$data[4,18] = 0.1371220078066022
$data[4,19] = 0.1371220078066022

# Row 7: MuSCs -> MuSCs
$data[5,0]  = "MuSCs"
$data[5,1]  = "Rln1"
$data[5,2]  = "Rxfp2"
$data[5,3]  = "MuSCs"
$data[5,4]  = 3
$data[5,5]  = 1
$data[5,6]  = 0.1012803333333333
$data[5,7]  = 0.303841
$data[5,8]  = 0.1387000293065843
$data[5,9]  = 0.1387000293065843
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.006558333333333333
$data[5,13] = 0.019675
$data[5,14] = 0.01137722542577134
$data[5,15] = 0.01137722542577134
$data[5,16] = 0.0006642301861111111
$data[5,17] = 0.005978071675
$data[5,18] = 0.001578021499982101
$data[5,19] = 0.001578021499982101

$ws.Range("A2:T7").Value = $data
